# Actualización automática 2025-09-08 13:10:08
#
# A new client, "BRAVO MANZABA MARIA CECILIA", was added to the
# "OFICINA-CATAECSA" asesor group (alphabetically right after
# "AVILA TORRES RAFAEL ALEJANDRO" and before "CARAVEDO PAZMIÑO  JAHAIRA
# PAMELA"). This inserts one new row into both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets, shifting every row below it down by one and
# bumping the "dimension"/footer counters accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": new row at 264, columns A:R
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows("264:264").Insert()
$ws1.Range("A264").Value = "OFICINA-CATAECSA"
$ws1.Range("B264").Value = "BRAVO MANZABA MARIA CECILIA"
$ws1.Range("C264:R264").Value = 0

# The footer row (previously row 316, now shifted to 317) holds text like
# "0 de 314" counting non-zero entries out of the total client count; the
# total must now read "315".
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(317, $col)
    $cell.Value = ($cell.Value2 -replace "314", "315")
}

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL": new row at 268, columns A:G
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows("268:268").Insert()
$ws2.Range("A268").Value = "OFICINA-CATAECSA"
$ws2.Range("B268").Value = "BRAVO MANZABA MARIA CECILIA"
$ws2.Range("C268:G268").Value = 0

# The footer row (previously row 320, now shifted to 321) is a plain sum
# of static values; since the new row is all zero it is unchanged, and the
# row-shift is handled automatically by Insert().
